$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - Spanish descriptive headers -> short codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of "de"/"del"/"la"/"los"/"el" particles to title case
$ws.Range("B15").Value = "Ocozocoautla De Espinosa"
$ws.Range("B17").Value = "San Juan De Sabinas"
$ws.Range("B19").Value = "Villa De Álvarez"
$ws.Range("A21").Value = "Ciudad De México"
$ws.Range("A35").Value = "Estado De México"
$ws.Range("B36").Value = "Ecatepec De Morelos"
$ws.Range("B37").Value = "Ixtapan De La Sal"
$ws.Range("B41").Value = "San Felipe Del Progreso"
$ws.Range("B45").Value = "Tlalnepantla De Baz"
$ws.Range("B56").Value = "Acapulco De Juárez"
$ws.Range("B57").Value = "Ayutla De Los Libres"
$ws.Range("B58").Value = "Buenavista De Cuéllar"
$ws.Range("B69").Value = "Jilotlán De Los Dolores"
$ws.Range("B71").Value = "La Manzanilla De La Paz"
$ws.Range("B75").Value = "Tizapán El Alto"
$ws.Range("B76").Value = "Unión De Tula"
$ws.Range("B79").Value = "Zapotlán El Grande"
$ws.Range("B90").Value = "Ixtlán Del Río"
$ws.Range("B98").Value = "Chalcatongo De Hidalgo"
$ws.Range("B100").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B101").Value = "Ixtlán De Juárez"
$ws.Range("B102").Value = "Putla Villa De Guerrero"
$ws.Range("B113").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B117").Value = "Huehuetlán El Chico"
$ws.Range("B118").Value = "Los Reyes De Juárez"
$ws.Range("B125").Value = "Tetela De Ocampo"
$ws.Range("B148").Value = "Ignacio De La Llave"

# Remove the trailing metadata/footer rows (previously rows 156-160),
# which shrinks the used dimension from A1:D160 to A1:D154
$ws.Range("A156:A160").EntireRow.Delete() | Out-Null

